# Update "想去人数" (F column) counts that increased by the time the data
# was regenerated. The same underlying events appear on both the "展览"
# sheet and the aggregated "全部类型" sheet, so each value must be updated
# in both places.

$wb = $excel.ActiveWorkbook

# Sheet "展览" (rows keyed by row number)
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F2").Value = 1870
$wsExhibit.Range("F5").Value = 174
$wsExhibit.Range("F6").Value = 2572
$wsExhibit.Range("F22").Value = 57
$wsExhibit.Range("F23").Value = 1644
$wsExhibit.Range("F25").Value = 401
$wsExhibit.Range("F28").Value = 299

# Sheet "全部类型" (same events, but offset by one row after row 2 because
# the "演出" sheet's single row is inserted as row 3)
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F2").Value = 1870
$wsAll.Range("F6").Value = 174
$wsAll.Range("F7").Value = 2572
$wsAll.Range("F23").Value = 57
$wsAll.Range("F24").Value = 1644
$wsAll.Range("F26").Value = 401
$wsAll.Range("F29").Value = 299
